# Updates cryptos list figures (price/volume columns D & E) and restores two
# swapped coin rows (28/29: BinanceUSD<->EthereumClassic, 42/43: PaxDollar<->WEMIXToken)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. A leading "'" forces text (prevents Excel from
# auto-parsing values like "1.00" / "27.625.06" as numbers); ClearFormats() afterwards
# strips the quote-prefix artifact and any implicit numeric style, so cells keep the
# same (unstyled) look as every other inline-string cell in the sheet.
$updates = [ordered]@{
    "D2" = "'27.625.06"
    "E2" = "'  +0.77%  "
    "D3" = "'1.636.01"
    "E3" = "'  -0.39%  "
    "E4" = "'  +0.14%  "
    "D5" = "'212.55"
    "E5" = "'  +0.31%  "
    "D6" = "'0.522"
    "E6" = "'  -1.28%  "
    "E7" = "'  +0.18%  "
    "D8" = "'22.94"
    "E8" = "'  -0.68%  "
    "E9" = "'  +0.46%  "
    "E10" = "'  -0.10%  "
    "E11" = "'  +0.62%  "
    "D12" = "'1.868.69"
    "E12" = "'  -0.33%  "
    "D13" = "'1.667.17"
    "E13" = "'  +1.65%  "
    "E14" = "'  -0.06%  "
    "E15" = "'  -1.81%  "
    "D16" = "'64.52"
    "E16" = "'  +0.17%  "
    "D17" = "'27.624.67"
    "E17" = "'  +0.86%  "
    "D18" = "'229.39"
    "E18" = "'  +0.07%  "
    "D19" = "'7.73"
    "E19" = "'  +1.59%  "
    "E20" = "'  +0.16%  "
    "E21" = "'  +0.13%  "
    "E22" = "'  -1.19%  "
    "D23" = "'10.02"
    "E23" = "'  +4.19%  "
    "E24" = "'  -2.86%  "
    "D25" = "'150.21"
    "E25" = "'  +2.04%  "
    "E26" = "'  -1.13%  "
    "E27" = "'  -1.53%  "
    "B28" = "'BinanceUSD"
    "C28" = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
    "D28" = "'1.00"
    "E28" = "'  +0.20%  "
    "B29" = "'EthereumClassic"
    "C29" = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D29" = "'15.62"
    "E29" = "'  +0.32%  "
    "E30" = "'  +0.17%  "
    "D31" = "'0.0485"
    "E31" = "'  -0.02%  "
    "D32" = "'3.29"
    "E32" = "'  +0.21%  "
    "D33" = "'1.453.23"
    "E33" = "'  +2.73%  "
    "D34" = "'3.11"
    "E34" = "'  -1.60%  "
    "E35" = "'  -0.97%  "
    "E36" = "'  +0.10%  "
    "D37" = "'0.564"
    "E37" = "'  -0.37%  "
    "D38" = "'0.874"
    "E38" = "'  -1.19%  "
    "D40" = "'0.895"
    "E40" = "'  +8.90%  "
    "D41" = "'69.74"
    "E41" = "'  +8.03%  "
    "B42" = "'PaxDollar"
    "C42" = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "D42" = "'1.00"
    "E42" = "'  +0.18%  "
    "B43" = "'WEMIXToken"
    "C43" = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D43" = "'1.02"
    "E43" = "'  -0.94%  "
    "D45" = "'2.46"
    "E45" = "'  +0.39%  "
    "E46" = "'  -0.14%  "
    "D47" = "'1.778.69"
    "E47" = "'  -0.33%  "
    "E48" = "'  +1.94%  "
    "D49" = "'86.21"
    "E49" = "'  -2.05%  "
    "E50" = "'  -1.22%  "
    "D51" = "'0.0986"
    "E51" = "'  -0.45%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $updates[$cellRef]
    $cell.ClearFormats()
}

